$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.719.65'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '1.742.00'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +1.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.82'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3835'
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3344'
$ws.Range("E8").Value = '  -2.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.37'
$ws.Range("E9").Value = '  -4.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.099'
$ws.Range("E10").Value = '  -4.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07162'
$ws.Range("E11").Value = '  -3.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.13'
$ws.Range("E13").Value = '  -3.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.099'
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").Value = '1.740.75'
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.949'
$ws.Range("E16").Value = '  -2.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001047'
$ws.Range("E17").Value = '  -3.08%  '
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.25'
$ws.Range("E20").Value = '  -5.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.58'
$ws.Range("E21").Value = '  -5.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.137'
$ws.Range("E22").Value = '  -4.57%  '
$ws.Range("D23").Value = '27.715.83'
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.51'
$ws.Range("E24").Value = '  -5.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.397'
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.66'
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.59'
$ws.Range("E27").Value = '  -5.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.254'
$ws.Range("E28").Value = '  -7.80%  '
$ws.Range("D29").Value = '1.940.54'
$ws.Range("E29").Value = '  -1.84%  '
$ws.Range("E30").Value = '  -12.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '128.40'
$ws.Range("E31").Value = '  -4.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.018'
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.739'
$ws.Range("E33").Value = '  -7.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08670'
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.93'
$ws.Range("E35").Value = '  -7.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.521'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6421'
$ws.Range("E37").Value = '  -6.91%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.066'
$ws.Range("E38").Value = '  -5.27%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02246'
$ws.Range("E39").Value = '  -7.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06026'
$ws.Range("E40").Value = '  -5.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2075'
$ws.Range("E41").Value = '  -5.51%  '
$ws.Range("E42").Value = '  -4.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.887'
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.46'
$ws.Range("E45").Value = '  -5.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.791'
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5944'
$ws.Range("E47").Value = '  -6.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.66'
$ws.Range("E48").Value = '  -5.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.966'
$ws.Range("E49").Value = '  -6.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.145'
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06922'
$ws.Range("E51").Value = '  -6.85%  '
